## Updated xlsform template to include display_name for choices
#
# 1. "choices" sheet: insert a new column C ("display_name") between
#    "name" (B) and "label::language" (old C, now D). The frozen pane /
#    split moves one column to the right along with it.
# 2. Bump zoom to 150% on both the "survey" and "choices" sheet views.
# 3. Leave "survey" as the active/selected sheet when done (switching to
#    "choices" to edit it would otherwise strand tabSelected there).

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- choices sheet: insert the new "display_name" column -----------------
$choices.Activate()
$choices.Columns.Item(3).Insert()
$choices.Cells.Item(1, 3).Value = "display_name"

# Recreate the frozen pane one column further right (was C/D split at
# D2, now D/E split at E2), then leave the selection on the new header
# cell, matching the saved view state.
$excel.ActiveWindow.FreezePanes = $false
[void]$choices.Range("E2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$choices.Range("C1").Select()

# 150% zoom on the choices view.
$excel.ActiveWindow.Zoom = 150

# --- survey sheet: zoom only, no structural change ------------------------
$survey.Activate()
[void]$survey.Range("A3").Select()
$excel.ActiveWindow.Zoom = 150
